# DSSP.xlsx: move from a single "img" column / 1 product to a 3-image-per-product
# layout (img1 / img2 / img3 in columns C/D/E) with 4 sample products.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (barcode/name already existed; img1/2/3 come later) ----
$ws.Range("A1").Value = "barcode"
$ws.Range("B1").Value = "name"

# ---- Product rows ------------------------------------------------------
$ws.Range("A2").Value = 8850006325636
$ws.Range("B2").Value = "KDR Colgate TOT ActiveFresh 150g"
$ws.Range("C2").Value = "hinh_anh_san_pham\KDR_Colgate_TOT_ActiveFresh_150g_1.jpg"
$ws.Range("D2").Value = "hinh_anh_san_pham\KDR_Colgate_TOT_ActiveFresh_150g_2.jpg"
$ws.Range("E2").Value = "hinh_anh_san_pham\KDR_Colgate_TOT_ActiveFresh_150g_3.jpg"

$ws.Range("A3").Value = 8850006327647
$ws.Range("B3").Value = "KDR Colgate CSPR Comp Protect 110g"
$ws.Range("C3").Value = "hinh_anh_san_pham\KDR_Colgate_CSPR_Comp_Protect_110g_1.jpg"
$ws.Range("D3").Value = "hinh_anh_san_pham\KDR_Colgate_CSPR_Comp_Protect_110g_2.jpg"
$ws.Range("E3").Value = "hinh_anh_san_pham\KDR_Colgate_CSPR_Comp_Protect_110g_3.jpg"

$ws.Range("A4").Value = 8850006331866
$ws.Range("B4").Value = "BCDR Colgate SlimSoft Charcoal 1PK"
$ws.Range("C4").Value = "hinh_anh_san_pham\BCDR_Colgate_SlimSoft_Charcoal_1PK_1.jpg"
$ws.Range("D4").Value = "hinh_anh_san_pham\BCDR_Colgate_SlimSoft_Charcoal_1PK_2.jpg"
$ws.Range("E4").Value = "hinh_anh_san_pham\BCDR_Colgate_SlimSoft_Charcoal_1PK_3.jpg"

$ws.Range("A5").Value = 8850006332030
$ws.Range("B5").Value = "BCDR Colgate 360 Char Spiral 2"
$ws.Range("C5").Value = "hinh_anh_san_pham\BCDR_Colgate_360_Char_Spiral_2_1.jpg"
$ws.Range("D5").Value = "hinh_anh_san_pham\BCDR_Colgate_360_Char_Spiral_2_2.jpg"
$ws.Range("E5").Value = "hinh_anh_san_pham\BCDR_Colgate_360_Char_Spiral_2_3.jpg"

# ---- Header for the 3 image columns (written last, like the script that
# labelled the already-populated columns) --------------------------------
$ws.Range("C1").Value = "img1"
$ws.Range("D1").Value = "img2"
$ws.Range("E1").Value = "img3"

# New image-path cells use the plain (unstyled) look already used by the
# original "img" column, i.e. no border/alignment formatting.
$ws.Range("C2:E5").Style = "Normal"

# ---- Row heights (row 2 grew slightly taller; row 4 slightly shorter) --
$ws.Rows.Item(2).RowHeight = 67.9
$ws.Rows.Item(3).RowHeight = 67.9
$ws.Rows.Item(4).RowHeight = 41.45
$ws.Rows.Item(5).RowHeight = 54.6
$ws.Rows.Item(6).RowHeight = 54.6

# ---- Column widths: C stays a "bestFit" image-path column; D/E are new
# image-path columns sized the same way; F:H keep the sheet's base width. --
$ws.Columns.Item(1).ColumnWidth = 32.28515625
$ws.Columns.Item(2).ColumnWidth = 28.28515625
$ws.Columns.Item(3).ColumnWidth = 64.28515625
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 62.42578125
$ws.Range("F1:H1").EntireColumn.ColumnWidth = 8.85546875

# ---- Selection, as last left by the editing session --------------------
$null = $ws.Range("C6").Select()
